$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text changes: the old single title string in K1 is replaced by two
#    strings - the bowling center name (moved to N1) and a lane label
#    (new, J2). A new styled-but-empty cell is also created at R2.
# ---------------------------------------------------------------------------
$ws.Range("K1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats - reuse K1's bold style
$ws.Range("N1").Value = "ABC Bowling Center"

$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").Value = "Lane  12"

$ws.Range("R2").PasteSpecial(-4122)

$ws.Range("K1").Clear()

# ---------------------------------------------------------------------------
# 2. Pictures: two of the five pictures on the sheet are removed, and the
#    remaining three are resized/repositioned (taller & narrower, shifted
#    one column to the left) to cover the enlarged print area. This must
#    happen *before* any column-width edits below, since shifting a
#    column's width moves the pixel position of every column after it.
# ---------------------------------------------------------------------------
$ws.Shapes.Item("Picture 4").Delete()
$ws.Shapes.Item("Picture 5").Delete()

# Picture 1 keeps its original top-left anchor (column B / row 3) but grows
# taller and a bit narrower.
$p1 = $ws.Shapes.Item("Picture 1")
$p1.Left = $ws.Columns.Item(2).Left
$p1.Top = $ws.Rows.Item(3).Top
$p1.Width = 450
$p1.Height = 1100

# Picture 2 shifts one column to the left (O -> N) and gets the same new size.
$p2 = $ws.Shapes.Item("Picture 2")
$p2.Left = $ws.Columns.Item(15).Left
$p2.Top = $ws.Rows.Item(3).Top
$p2.Width = 450
$p2.Height = 1100

# Picture 3 shifts one column to the left (O -> N) too; size is unchanged.
$p3 = $ws.Shapes.Item("Picture 3")
$p3.Left = $ws.Columns.Item(14).Left
$p3.Top = $ws.Rows.Item(11).Top
$p3.Width = 65
$p3.Height = 350

# ---------------------------------------------------------------------------
# 3. Column width tweaks that accompany the new layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 13.928566666666667
$ws.Columns.Item(17).ColumnWidth = 13.452376666666666
